# historico.xlsx — "Add files via upload" re-commit.
#
# The ranking column D (header "posicao") on Planilha1 is renumbered so it
# runs sequentially 1..41 alongside the data, which is already sorted by
# column C (value) descending. Rows 2-13 already hold D=1..12; rows 14-42
# get updated here to D=13..41.
#
# The sheet's active selection is also left on a single cell (C4) instead
# of the whole table range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

for ($row = 14; $row -le 42; $row++) {
    $ws.Cells.Item($row, 4).Value = $row - 1
}

[void]$ws.Range("C4").Select()
